# Added missing payment from Ganesh Kumar
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ganesh Kumar is row 8 (A8=6, B8="Ganesh Kumar"). He had a second payment
# of 400 that was missing from column E - add it in. All the dependent
# totals/running-balance formulas (E35, P35, E40:E52, ...) recompute
# automatically from this single input change.
$ws.Range("E8").Value = 400

# Restore the view/selection state captured with the edit: window scrolled
# so A10 is the top-left visible cell, with R23 as the active selection.
$excel.Goto($ws.Range("A10"), $true)
$ws.Range("R23").Select()

$wb.Save()
